$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.115.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.28%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.411.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.87%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'578.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.05%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'7.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +6.95%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +4.31%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.991.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.90%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.59%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +5.06%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.406.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.96%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'25.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.71%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.151.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.26%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +6.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.27%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'9.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.50%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'389.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +9.98%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +2.09%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.544.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.80%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +15.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'71.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.14%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.65%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.33%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'8.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.52%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +3.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +2.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'RenzoRestakedETH"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'3.441.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.84%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'USDe"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.03%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'23.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.51%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.04%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.60%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'163.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.85%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +13.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.787"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.75%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.21%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'4.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'25.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.51%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'41.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.29%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'6.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.54%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'23.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.85%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.372.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +8.34%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0266"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +5.51%  "
$ws.Range("E51").Style = "Normal"
